# Commiting changes of email reports
# Adds a "DB Type" column (new column A) and a second query row (MYSQL)
# to the "Queries" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Queries")

# --- Structural changes -------------------------------------------------
# Insert a new column before the existing column A; this shifts the old
# A:G content to B:H.
$ws.Columns.Item(1).Insert()

# Insert a new row before the (soon to exist) row 3, giving us a blank
# row 3 to populate with the MYSQL query entry.
$ws.Rows.Item(3).Insert()

# --- New column A: "DB Type" --------------------------------------------
$ws.Range("A1").Value = "DB Type"
$ws.Range("A2").Value = "MS SQL"
$ws.Range("A3").Value = "MYSQL"
$ws.Range("A1:A3").VerticalAlignment = -4160

# --- New row 3: mirrors row 2 in columns B:D, new data in E:H ----------
$ws.Range("B3").Value = $ws.Range("B2").Text
$ws.Range("C3").Value = $ws.Range("C2").Text
$ws.Range("D3").Value = $ws.Range("D2").Text
$ws.Range("E3").Value = "01-06-2021 00:00:00"
$ws.Range("F3").Value = "07-06-2021 00:00:00"
$ws.Range("G3").Value = "Select * from AGT_Agent_TimeTrack;"
$ws.Range("H3").Value = $ws.Range("H2").Text

# Re-apply the date/query-cell formatting (quote-prefixed, top-aligned)
# that the plain .Value assignment above reset on E3/F3.
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)

# Row 3 needs a tall, wrapped row to hold the long SQL text.
$ws.Rows.Item(3).RowHeight = 406

# --- Column sizing for the two new-ish columns ---------------------------
$ws.Columns.Item(1).ColumnWidth = 7.9    # -> stored width ~8.73 ("DB Type")
$ws.Columns.Item(2).ColumnWidth = 12.95  # -> stored width ~13.82 ("IVR"/"MS SQL"/"MYSQL")

# --- Selection / view ----------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("G3").Select() | Out-Null

Write-Output "done"
